# The commit swaps the contents of the two theme parts in this deck:
#   ppt/theme/theme1.xml (the slide master's theme, currently "Integral")
#   ppt/theme/theme2.xml (the notes master's theme, currently "Office Theme")
# become each other's colour scheme. We reproduce that through the
# PowerPoint object model by pushing the "Office Theme" palette onto the
# presentation's (slide-master) ThemeColorScheme, which is the theme that
# PowerPoint's automation surface exposes for editing.

$p = $ppt.ActivePresentation

# VBA's RGB() packs bytes as 0x00BBGGRR; build the same integer by hand
# since the PowerShell host here has no RGB() builtin.
function ColorInt([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

# Target palette: the "Office Theme" colours that currently live in
# ppt/theme/theme2.xml, in ThemeColorScheme order
# (1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink).
$officeTheme = @(
    @(0x00, 0x00, 0x00), # dk1
    @(0xFF, 0xFF, 0xFF), # lt1
    @(0x44, 0x54, 0x6A), # dk2
    @(0xE7, 0xE6, 0xE6), # lt2
    @(0x5B, 0x9B, 0xD5), # accent1
    @(0xED, 0x7D, 0x31), # accent2
    @(0xA5, 0xA5, 0xA5), # accent3
    @(0xFF, 0xC0, 0x00), # accent4
    @(0x44, 0x72, 0xC4), # accent5
    @(0x70, 0xAD, 0x47), # accent6
    @(0x05, 0x63, 0xC1), # hlink
    @(0x95, 0x4F, 0x72)  # folHlink
)

$master = $p.Designs.Item(1).SlideMaster
$tcs = $master.Theme.ThemeColorScheme

for ($i = 0; $i -lt $officeTheme.Count; $i++) {
    $rgb = $officeTheme[$i]
    $tcs.Item($i + 1).RGB = ColorInt $rgb[0] $rgb[1] $rgb[2]
}
